# Apply updated cryptocurrency price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.482.00'
$ws.Range('E2').Value = '  +9.66%  '
$ws.Range('D3').Value = '3.235.04'
$ws.Range('E3').Value = '  +4.44%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '401.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.49'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.36%  '
$ws.Range('E7').Value = '  +3.30%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.628'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +7.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.71'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +7.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0903'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.52%  '
$ws.Range('D13').Value = '3.746.73'
$ws.Range('E13').Value = '  +4.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.25'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.12'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.53%  '
$ws.Range('E16').Value = '  +7.89%  '
$ws.Range('D17').Value = '3.254.12'
$ws.Range('E17').Value = '  +5.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.73'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.86%  '
$ws.Range('D19').Value = '56.279.27'
$ws.Range('E19').Value = '  +9.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.45'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000104'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +7.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.21'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '295.15'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +11.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.75'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.87%  '
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.25'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.27'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.68%  '
$ws.Range('E28').Value = '  +2.92%  '
$ws.Range('E29').Value = '  +2.87%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.39'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +10.40%  '
$ws.Range('E33').Value = '  +5.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '36.85'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.11'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.50'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('E37').Value = '  +6.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.14'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +24.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '136.74'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.94'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.05'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +10.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.25'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.26%  '
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('E45').Value = '  -5.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.68'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.13'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +41.63%  '
$ws.Range('D48').Value = '2.159.97'
$ws.Range('E48').Value = '  +5.18%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.43'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0360'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +9.74%  '
